$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 9321928.0
$ws.Range("B7").Value = 7990280.0
$ws.Range("A8").Value = 9192224.0
$ws.Range("B8").Value = 7879120.0
